$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '45.176.49'
$ws.Range("E2").Value = '  +3.28%  '

$ws.Range("D3").Value = '2.364.12'
$ws.Range("E3").Value = '  +1.37%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").Value = '''311.80'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.50%  '

$ws.Range("D6").Value = '''108.58'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.25%  '

$ws.Range("D7").Value = '''0.631'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.11%  '

$ws.Range("E8").Value = '  -0.14%  '

$ws.Range("D9").Value = '''0.610'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.62%  '

$ws.Range("D10").Value = '''40.83'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.13%  '

$ws.Range("D11").Value = '''0.0915'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.22%  '

$ws.Range("D12").Value = '''8.45'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.09%  '

$ws.Range("E13").Value = '  +1.08%  '

$ws.Range("E14").Value = '  -3.53%  '

$ws.Range("D15").Value = '2.724.22'
$ws.Range("E15").Value = '  +1.41%  '

$ws.Range("D16").Value = '''15.21'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.69%  '

$ws.Range("D17").Value = '2.360.98'
$ws.Range("E17").Value = '  +1.39%  '

$ws.Range("D18").Value = '45.127.57'
$ws.Range("E18").Value = '  +3.24%  '

$ws.Range("D19").Value = '''14.38'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +10.60%  '

$ws.Range("E20").Value = '  -0.89%  '

$ws.Range("D21").Value = '''7.19'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.01%  '

$ws.Range("D22").Value = '''73.09'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.69%  '

$ws.Range("D23").Value = '''3.49'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.30%  '

$ws.Range("D24").Value = '''258.75'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.64%  '

$ws.Range("D25").Value = '''2.31'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.13%  '

$ws.Range("E26").Value = '  -0.20%  '

$ws.Range("D27").Value = '''11.06'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.76%  '

$ws.Range("D28").Value = '''7.21'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -5.42%  '

$ws.Range("E29").Value = '  +0.47%  '

$ws.Range("D30").Value = '''0.0971'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +7.82%  '

$ws.Range("D31").Value = '''22.36'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.80%  '

$ws.Range("D32").Value = '''37.12'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -5.99%  '

$ws.Range("D33").Value = '''167.91'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.13%  '

$ws.Range("D34").Value = '''3.01'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +5.45%  '

$ws.Range("E35").Value = '  -1.30%  '

$ws.Range("D36").Value = '''0.117'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.63%  '

$ws.Range("D37").Value = '''4.67'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.74%  '

$ws.Range("D38").Value = '''3.96'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.62%  '

$ws.Range("D39").Value = '''0.0353'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.36%  '

$ws.Range("D40").Value = '''2.88'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.92%  '

$ws.Range("D41").Value = '''1.78'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.46%  '

$ws.Range("D42").Value = '''99.44'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.42%  '

$ws.Range("B43").Value = 'MultiversX'
$ws.Range("C43").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D43").Value = '''69.38'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.25%  '

$ws.Range("B44").Value = 'Algorand'
$ws.Range("C44").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D44").Value = '''0.229'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.93%  '

$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").Value = '1.868.48'
$ws.Range("E45").Value = '  +12.64%  '

$ws.Range("D46").Value = '''12.84'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -6.41%  '

$ws.Range("E47").Value = '  -0.45%  '

$ws.Range("D48").Value = '''83.39'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +9.87%  '

$ws.Range("D49").Value = '''5.67'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +7.29%  '

$ws.Range("D50").Value = '''9.16'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.83%  '

$ws.Range("D51").Value = '''110.11'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.24%  '

Write-Output "done"
